$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = "A1"
$ws.Range("B2").Value = "A2"
$ws.Range("B3").Value = "A3"
$ws.Range("B4").Value = "A4"
